# Apply crypto price/volume updates as described by the commit diff.
# Some "Price" (column D) values, although stored as text in the source
# workbook, look like plain numbers (e.g. "3.45"). Assigning such a
# string directly to Range.Value causes Excel to auto-convert it into a
# numeric cell, which would not match the original text-cell layout.
# Set-TextValue forces the value to be written/kept as text, then
# restores the cell to the default "Normal" style so no stray
# number-format styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "73.270.33"
$ws.Range("E2").Value = "  +0.57%  "
Set-TextValue "D3" "3.964.76"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "609.88"
$ws.Range("E5").Value = "  +8.80%  "
Set-TextValue "D6" "168.51"
$ws.Range("E6").Value = "  +11.12%  "
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  +7.00%  "
Set-TextValue "D11" "56.01"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("E12").Value = "  +2.29%  "
Set-TextValue "D13" "11.23"
$ws.Range("E13").Value = "  +1.94%  "
Set-TextValue "D14" "4.609.09"
$ws.Range("E14").Value = "  -2.04%  "
Set-TextValue "D15" "3.969.71"
$ws.Range("E15").Value = "  -2.28%  "
Set-TextValue "D16" "14.15"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").Value = "  +1.91%  "
Set-TextValue "D18" "20.41"
$ws.Range("E18").Value = "  -2.06%  "
Set-TextValue "D19" "73.155.14"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  -1.07%  "
Set-TextValue "D21" "438.38"
$ws.Range("E21").Value = "  -1.83%  "
Set-TextValue "D22" "4.82"
$ws.Range("E22").Value = "  +9.19%  "
Set-TextValue "D23" "95.50"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("E24").Value = "  -5.47%  "
Set-TextValue "D25" "14.18"
$ws.Range("E25").Value = "  -4.45%  "
Set-TextValue "D26" "4.08"
$ws.Range("E26").Value = "  -7.36%  "
$ws.Range("E27").Value = "  -3.23%  "
Set-TextValue "D28" "5.95"
$ws.Range("E28").Value = "  +0.08%  "
Set-TextValue "D29" "10.45"
Set-TextValue "D30" "35.99"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("E31").Value = "  -2.38%  "
Set-TextValue "D32" "13.84"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +17.73%  "
$ws.Range("E34").Value = "  -4.45%  "
Set-TextValue "D35" "47.68"
$ws.Range("E35").Value = "  -3.36%  "
Set-TextValue "D36" "70.26"
$ws.Range("E36").Value = "  +4.17%  "
Set-TextValue "D37" "648.06"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("E38").Value = "  -5.37%  "
Set-TextValue "D39" "3.39"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  -0.10%  "
Set-TextValue "D41" "0.146"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  +0.10%  "
Set-TextValue "D43" "0.0483"
$ws.Range("E43").Value = "  -2.89%  "
Set-TextValue "D44" "3.20"
$ws.Range("E44").Value = "  -5.96%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D45" "10.50"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D46" "3.12"
$ws.Range("E46").Value = "  +32.39%  "
Set-TextValue "D47" "0.000302"
$ws.Range("E47").Value = "  +8.76%  "
$ws.Range("E48").Value = "  -2.86%  "
Set-TextValue "D49" "3.45"
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("E50").Value = "  -6.03%  "
Set-TextValue "D51" "2.98"
$ws.Range("E51").Value = "  -4.73%  "
